# Remove user authentication module
# (Appends the latest sensor reading (row 45) to each of the four
#  sensor-log worksheets, mirroring the existing row layout.)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW35-FE-LIFTER ---------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("A45").Value = "2025-03-06 04:42:06"
$ws.Range("B45").Value = "0x01,0x90 "
$ws.Range("C45").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Range("D45").Value = "0x01,0x90,"
$ws.Range("E45").Value = "0x d"
$ws.Range("F45").Value = 400
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "568631262647113770877196"
$ws.Range("H45").Value = 400
$ws.Range("I45").Value = 13

# --- Sheet 2: ROW35-MID-LIFTER --------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("A45").Value = "2025-03-06 04:29:35"
$ws.Range("B45").Value = "0x01,0x90 "
$ws.Range("C45").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Range("D45").Value = "0x01,0x90,"
$ws.Range("E45").Value = "0x e"
$ws.Range("F45").Value = 400
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "568631262647113770942732"
$ws.Range("H45").Value = 400
$ws.Range("I45").Value = 14

# --- Sheet 3: ROW02-FE-LIFTER ---------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("A45").Value = "2025-03-06 04:51:45"
$ws.Range("B45").Value = "0x01,0x90 "
$ws.Range("C45").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
$ws.Range("D45").Value = "0x01,0x90,"
$ws.Range("E45").Value = "0xff"
$ws.Range("F45").Value = 400
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "568631262647113769959692"
$ws.Range("H45").Value = 400
$ws.Range("I45").Value = 255

# --- Sheet 4: ROW02-MID-LIFTER --------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("A45").Value = "2025-03-06 04:41:15"
$ws.Range("B45").Value = "0x01,0x90 "
$ws.Range("C45").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D45").Value = "0x01,0x90,"
$ws.Range("E45").Value = "0x 3"
$ws.Range("F45").Value = 400
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "568631262647113769959692"
$ws.Range("H45").Value = 400
$ws.Range("I45").Value = 3
